$wb = $excel.ActiveWorkbook

# ==== Sheet 1: strategy_id-0 ====
$ws1 = $wb.Worksheets.Item(1)
$arr = New-Object 'object[,]' 1,4
$arr[0,0] = "variable_trajectory_group"
$arr[0,1] = "normalize_group"
$arr[0,2] = "trajgroup_no_vary_q"
$arr[0,3] = "uniform_scaling_q"
$ws1.Range("C1:F1").Value = $arr
$arr = New-Object 'object[,]' 1,2
$arr[0,0] = "max_55"
$arr[0,1] = "min_55"
$ws1.Range("H1:I1").Value = $arr
$ws1.Range("J1").Copy($ws1.Range("AT1:BM1"))
$arr = New-Object 'object[,]' 1,20
$arr[0,0] = 36
$arr[0,1] = 37
$arr[0,2] = 38
$arr[0,3] = 39
$arr[0,4] = 40
$arr[0,5] = 41
$arr[0,6] = 42
$arr[0,7] = 43
$arr[0,8] = 44
$arr[0,9] = 45
$arr[0,10] = 46
$arr[0,11] = 47
$arr[0,12] = 48
$arr[0,13] = 49
$arr[0,14] = 50
$arr[0,15] = 51
$arr[0,16] = 52
$arr[0,17] = 53
$arr[0,18] = 54
$arr[0,19] = 55
$ws1.Range("AT1:BM1").Value = $arr
$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 6257.313670200001
$arr[0,1] = 6536.446673800001
$arr[0,2] = 6815.579677400001
$arr[0,3] = 7094.712681000001
$arr[0,4] = 7571.840713400001
$arr[0,5] = 8048.968745800003
$arr[0,6] = 8526.096778200004
$arr[0,7] = 9003.224810600004
$arr[0,8] = 9480.352843000004
$arr[0,9] = 9999.181414400004
$ws1.Range("L2:U2").Value = $arr
$arr = New-Object 'object[,]' 1,14
$arr[0,0] = 16402.39476400002
$arr[0,1] = 17083.92902600002
$arr[0,2] = 17765.46328800002
$arr[0,3] = 18446.99755000002
$arr[0,4] = 19206.92426600002
$arr[0,5] = 19966.85098200002
$arr[0,6] = 20726.77769800002
$arr[0,7] = 21486.70441400002
$arr[0,8] = 22246.63113000002
$arr[0,9] = 23076.36723400002
$arr[0,10] = 23906.10333800003
$arr[0,11] = 24735.83944200003
$arr[0,12] = 25565.57554600003
$arr[0,13] = 26395.31165000003
$ws1.Range("AF2:AS2").Value = $arr
$arr = New-Object 'object[,]' 1,20
$arr[0,0] = 27273.40450150054
$arr[0,1] = 28162.23718332734
$arr[0,2] = 29060.93273098747
$arr[0,3] = 29968.56370866613
$arr[0,4] = 30884.15306034566
$arr[0,5] = 31806.67516575942
$arr[0,6] = 32735.05710360175
$arr[0,7] = 33668.18012351037
$arr[0,8] = 34604.88132738751
$arr[0,9] = 35543.95555963519
$arr[0,10] = 36484.15750485686
$arr[0,11] = 37424.20399052605
$arr[0,12] = 38362.77649104803
$arr[0,13] = 39298.52382855237
$arr[0,14] = 40230.06506465726
$arr[0,15] = 41155.99257634703
$arr[0,16] = 42074.8753080102
$arr[0,17] = 42985.26219060671
$arr[0,18] = 43885.68571787133
$arr[0,19] = 44774.66566842715
$ws1.Range("AT2:BM2").Value = $arr
$arr = New-Object 'object[,]' 1,20
$arr[0,0] = 328726000
$arr[0,1] = 328726000
$arr[0,2] = 328726000
$arr[0,3] = 328726000
$arr[0,4] = 328726000
$arr[0,5] = 328726000
$arr[0,6] = 328726000
$arr[0,7] = 328726000
$arr[0,8] = 328726000
$arr[0,9] = 328726000
$arr[0,10] = 328726000
$arr[0,11] = 328726000
$arr[0,12] = 328726000
$arr[0,13] = 328726000
$arr[0,14] = 328726000
$arr[0,15] = 328726000
$arr[0,16] = 328726000
$arr[0,17] = 328726000
$arr[0,18] = 328726000
$arr[0,19] = 328726000
$ws1.Range("AT3:BM3").Value = $arr
$arr = New-Object 'object[,]' 1,20
$arr[0,0] = 1
$arr[0,1] = 1
$arr[0,2] = 1
$arr[0,3] = 1
$arr[0,4] = 1
$arr[0,5] = 1
$arr[0,6] = 1
$arr[0,7] = 1
$arr[0,8] = 1
$arr[0,9] = 1
$arr[0,10] = 1
$arr[0,11] = 1
$arr[0,12] = 1
$arr[0,13] = 1
$arr[0,14] = 1
$arr[0,15] = 1
$arr[0,16] = 1
$arr[0,17] = 1
$arr[0,18] = 1
$arr[0,19] = 1
$ws1.Range("AT4:BM4").Value = $arr
$arr = New-Object 'object[,]' 1,20
$arr[0,0] = -0.1
$arr[0,1] = -0.1
$arr[0,2] = -0.1
$arr[0,3] = -0.1
$arr[0,4] = -0.1
$arr[0,5] = -0.1
$arr[0,6] = -0.1
$arr[0,7] = -0.1
$arr[0,8] = -0.1
$arr[0,9] = -0.1
$arr[0,10] = -0.1
$arr[0,11] = -0.1
$arr[0,12] = -0.1
$arr[0,13] = -0.1
$arr[0,14] = -0.1
$arr[0,15] = -0.1
$arr[0,16] = -0.1
$arr[0,17] = -0.1
$arr[0,18] = -0.1
$arr[0,19] = -0.1
$ws1.Range("AT5:BM5").Value = $arr
$ws1.Range("C6").Value = 13
$arr = New-Object 'object[,]' 1,36
$arr[0,0] = 0.31
$arr[0,1] = 0.31
$arr[0,2] = 0.31
$arr[0,3] = 0.31
$arr[0,4] = 0.31
$arr[0,5] = 0.31
$arr[0,6] = 0.31
$arr[0,7] = 0.31
$arr[0,8] = 0.31
$arr[0,9] = 0.31
$arr[0,10] = 0.31
$arr[0,11] = 0.31
$arr[0,12] = 0.31
$arr[0,13] = 0.31
$arr[0,14] = 0.31
$arr[0,15] = 0.31
$arr[0,16] = 0.31
$arr[0,17] = 0.31
$arr[0,18] = 0.31
$arr[0,19] = 0.31
$arr[0,20] = 0.31
$arr[0,21] = 0.31
$arr[0,22] = 0.31
$arr[0,23] = 0.31
$arr[0,24] = 0.31
$arr[0,25] = 0.31
$arr[0,26] = 0.31
$arr[0,27] = 0.31
$arr[0,28] = 0.31
$arr[0,29] = 0.31
$arr[0,30] = 0.31
$arr[0,31] = 0.31
$arr[0,32] = 0.31
$arr[0,33] = 0.31
$arr[0,34] = 0.31
$arr[0,35] = 0.31
$ws1.Range("J6:AS6").Value = $arr
$arr = New-Object 'object[,]' 1,20
$arr[0,0] = 0.31
$arr[0,1] = 0.31
$arr[0,2] = 0.31
$arr[0,3] = 0.31
$arr[0,4] = 0.31
$arr[0,5] = 0.31
$arr[0,6] = 0.31
$arr[0,7] = 0.31
$arr[0,8] = 0.31
$arr[0,9] = 0.31
$arr[0,10] = 0.31
$arr[0,11] = 0.31
$arr[0,12] = 0.31
$arr[0,13] = 0.31
$arr[0,14] = 0.31
$arr[0,15] = 0.31
$arr[0,16] = 0.31
$arr[0,17] = 0.31
$arr[0,18] = 0.31
$arr[0,19] = 0.31
$ws1.Range("AT6:BM6").Value = $arr
$arr = New-Object 'object[,]' 1,20
$arr[0,0] = -999
$arr[0,1] = -999
$arr[0,2] = -999
$arr[0,3] = -999
$arr[0,4] = -999
$arr[0,5] = -999
$arr[0,6] = -999
$arr[0,7] = -999
$arr[0,8] = -999
$arr[0,9] = -999
$arr[0,10] = -999
$arr[0,11] = -999
$arr[0,12] = -999
$arr[0,13] = -999
$arr[0,14] = -999
$arr[0,15] = -999
$arr[0,16] = -999
$arr[0,17] = -999
$arr[0,18] = -999
$arr[0,19] = -999
$ws1.Range("AT7:BM7").Value = $arr
$arr = New-Object 'object[,]' 1,20
$arr[0,0] = -999
$arr[0,1] = -999
$arr[0,2] = -999
$arr[0,3] = -999
$arr[0,4] = -999
$arr[0,5] = -999
$arr[0,6] = -999
$arr[0,7] = -999
$arr[0,8] = -999
$arr[0,9] = -999
$arr[0,10] = -999
$arr[0,11] = -999
$arr[0,12] = -999
$arr[0,13] = -999
$arr[0,14] = -999
$arr[0,15] = -999
$arr[0,16] = -999
$arr[0,17] = -999
$arr[0,18] = -999
$arr[0,19] = -999
$ws1.Range("AT8:BM8").Value = $arr
$arr = New-Object 'object[,]' 1,20
$arr[0,0] = -999
$arr[0,1] = -999
$arr[0,2] = -999
$arr[0,3] = -999
$arr[0,4] = -999
$arr[0,5] = -999
$arr[0,6] = -999
$arr[0,7] = -999
$arr[0,8] = -999
$arr[0,9] = -999
$arr[0,10] = -999
$arr[0,11] = -999
$arr[0,12] = -999
$arr[0,13] = -999
$arr[0,14] = -999
$arr[0,15] = -999
$arr[0,16] = -999
$arr[0,17] = -999
$arr[0,18] = -999
$arr[0,19] = -999
$ws1.Range("AT9:BM9").Value = $arr
$arr = New-Object 'object[,]' 1,20
$arr[0,0] = 4.698333333333333
$arr[0,1] = 4.698333333333333
$arr[0,2] = 4.698333333333333
$arr[0,3] = 4.698333333333333
$arr[0,4] = 4.698333333333333
$arr[0,5] = 4.698333333333333
$arr[0,6] = 4.698333333333333
$arr[0,7] = 4.698333333333333
$arr[0,8] = 4.698333333333333
$arr[0,9] = 4.698333333333333
$arr[0,10] = 4.698333333333333
$arr[0,11] = 4.698333333333333
$arr[0,12] = 4.698333333333333
$arr[0,13] = 4.698333333333333
$arr[0,14] = 4.698333333333333
$arr[0,15] = 4.698333333333333
$arr[0,16] = 4.698333333333333
$arr[0,17] = 4.698333333333333
$arr[0,18] = 4.698333333333333
$arr[0,19] = 4.698333333333333
$ws1.Range("AT10:BM10").Value = $arr
$arr = New-Object 'object[,]' 1,34
$arr[0,0] = 888881381.2559999
$arr[0,1] = 892338114.0951002
$arr[0,2] = 895386227.1516801
$arr[0,3] = 898024053.4949002
$arr[0,4] = 906856770.0980221
$arr[0,5] = 910428206.9478642
$arr[0,6] = 913999643.7977062
$arr[0,7] = 917571080.6475483
$arr[0,8] = 921142517.4973904
$arr[0,9] = 923271168.4792284
$arr[0,10] = 925399819.4610667
$arr[0,11] = 927528470.4429048
$arr[0,12] = 929657121.4247428
$arr[0,13] = 931785772.4065808
$arr[0,14] = 932777441.0087162
$arr[0,15] = 933769109.6108512
$arr[0,16] = 934760778.2129862
$arr[0,17] = 935752446.8151214
$arr[0,18] = 936744115.4172565
$arr[0,19] = 936499448.5462054
$arr[0,20] = 936254781.6751541
$arr[0,21] = 936010114.8041033
$arr[0,22] = 935765447.9330523
$arr[0,23] = 935520781.0620012
$arr[0,24] = 934079056.6569703
$arr[0,25] = 932637332.2519392
$arr[0,26] = 931195607.8469083
$arr[0,27] = 929753883.4418772
$arr[0,28] = 928312159.0368462
$arr[0,29] = 925471525.2026653
$arr[0,30] = 922630891.3684845
$arr[0,31] = 919790257.5343034
$arr[0,32] = 916949623.7001225
$arr[0,33] = 914108989.8659414
$ws1.Range("L11:AS11").Value = $arr
$arr = New-Object 'object[,]' 1,20
$arr[0,0] = 910528891.5603858
$arr[0,1] = 906696344.2374226
$arr[0,2] = 902614158.694883
$arr[0,3] = 898285381.6000795
$arr[0,4] = 893713289.584669
$arr[0,5] = 888901383.1606724
$arr[0,6] = 883853380.4504673
$arr[0,7] = 878573210.7267025
$arr[0,8] = 873065007.7607515
$arr[0,9] = 867333102.9806495
$arr[0,10] = 861382018.4414332
$arr[0,11] = 855216459.6125952
$arr[0,12] = 848841307.9888465
$arr[0,13] = 842261613.5317464
$arr[0,14] = 835482586.9509113
$arr[0,15] = 828509591.834545
$arr[0,16] = 821348136.6399161
$arr[0,17] = 814003866.5551914
$arr[0,18] = 806482555.2447132
$arr[0,19] = 798790096.4903572
$ws1.Range("AT11:BM11").Value = $arr
$ws1.Range("Q12").Value = 507134002.4521359
$arr = New-Object 'object[,]' 1,25
$arr[0,0] = 551864940.7207718
$arr[0,1] = 563101060.9389338
$arr[0,2] = 574337181.1570959
$arr[0,3] = 585573301.3752576
$arr[0,4] = 596809421.5934196
$arr[0,5] = 608186310.7912844
$arr[0,6] = 619563199.9891495
$arr[0,7] = 630940089.1870143
$arr[0,8] = 642316978.3848792
$arr[0,9] = 653693867.5827446
$arr[0,10] = 665027848.8537955
$arr[0,11] = 676361830.1248465
$arr[0,12] = 687695811.3958975
$arr[0,13] = 699029792.6669484
$arr[0,14] = 710363773.9379994
$arr[0,15] = 721543059.9430304
$arr[0,16] = 732722345.9480615
$arr[0,17] = 743901631.9530923
$arr[0,18] = 755080917.9581234
$arr[0,19] = 766260203.9631542
$arr[0,20] = 776945460.7973353
$arr[0,21] = 787630717.6315161
$arr[0,22] = 798315974.4656971
$arr[0,23] = 809001231.2998779
$arr[0,24] = 819686488.134059
$ws1.Range("U12:AS12").Value = $arr
$arr = New-Object 'object[,]' 1,20
$arr[0,0] = 829543972.0022355
$arr[0,1] = 839197415.0471984
$arr[0,2] = 848636337.648142
$arr[0,3] = 857850394.4485067
$arr[0,4] = 866829391.9584272
$arr[0,5] = 875563306.1631435
$arr[0,6] = 884042300.0841366
$arr[0,7] = 892256741.2407371
$arr[0,8] = 900197218.9608759
$arr[0,9] = 907854561.4905095
$arr[0,10] = 915219852.8521475
$arr[0,11] = 922284449.403839
$arr[0,12] = 929039996.0509217
$arr[0,13] = 935478442.0638708
$arr[0,14] = 941592056.456735
$arr[0,15] = 947373442.8818125
$arr[0,16] = 952815553.99755
$arr[0,17] = 957911705.2680538
$arr[0,18] = 962655588.1541128
$arr[0,19] = 967041282.6572597
$ws1.Range("AT12:BM12").Value = $arr

# ==== Sheet 2: strategy_id-5004 ====
$ws2 = $wb.Worksheets.Item(2)
$arr = New-Object 'object[,]' 1,4
$arr[0,0] = "variable_trajectory_group"
$arr[0,1] = "normalize_group"
$arr[0,2] = "trajgroup_no_vary_q"
$arr[0,3] = "uniform_scaling_q"
$ws2.Range("C1:F1").Value = $arr
$arr = New-Object 'object[,]' 1,2
$arr[0,0] = "max_55"
$arr[0,1] = "min_55"
$ws2.Range("H1:I1").Value = $arr
$ws2.Range("J1").Copy($ws2.Range("AT1:BM1"))
$arr = New-Object 'object[,]' 1,20
$arr[0,0] = 36
$arr[0,1] = 37
$arr[0,2] = 38
$arr[0,3] = 39
$arr[0,4] = 40
$arr[0,5] = 41
$arr[0,6] = 42
$arr[0,7] = 43
$arr[0,8] = 44
$arr[0,9] = 45
$arr[0,10] = 46
$arr[0,11] = 47
$arr[0,12] = 48
$arr[0,13] = 49
$arr[0,14] = 50
$arr[0,15] = 51
$arr[0,16] = 52
$arr[0,17] = 53
$arr[0,18] = 54
$arr[0,19] = 55
$ws2.Range("AT1:BM1").Value = $arr
$ws2.Range("C2").Value = 13
$arr = New-Object 'object[,]' 1,36
$arr[0,0] = 0.31
$arr[0,1] = 0.31
$arr[0,2] = 0.31
$arr[0,3] = 0.31
$arr[0,4] = 0.31
$arr[0,5] = 0.31
$arr[0,6] = 0.31
$arr[0,7] = 0.31
$arr[0,8] = 0.31
$arr[0,9] = 0.31
$arr[0,10] = 0.31
$arr[0,11] = 0.31
$arr[0,12] = 0.31
$arr[0,13] = 0.31
$arr[0,14] = 0.31
$arr[0,15] = 0.31
$arr[0,16] = 0.31
$arr[0,17] = 0.31
$arr[0,18] = 0.31
$arr[0,19] = 0.31
$arr[0,20] = 0.31
$arr[0,21] = 0.31
$arr[0,22] = 0.3100000000000001
$arr[0,23] = 0.31
$arr[0,24] = 0.31
$arr[0,25] = 0.3100000000000001
$arr[0,26] = 0.3099999999999999
$arr[0,27] = 0.31
$arr[0,28] = 0.31
$arr[0,29] = 0.3099999999999999
$arr[0,30] = 0.31
$arr[0,31] = 0.31
$arr[0,32] = 0.31
$arr[0,33] = 0.31
$arr[0,34] = 0.31
$arr[0,35] = 0.31
$ws2.Range("J2:AS2").Value = $arr
$arr = New-Object 'object[,]' 1,20
$arr[0,0] = 0.31
$arr[0,1] = 0.31
$arr[0,2] = 0.31
$arr[0,3] = 0.31
$arr[0,4] = 0.31
$arr[0,5] = 0.31
$arr[0,6] = 0.31
$arr[0,7] = 0.31
$arr[0,8] = 0.31
$arr[0,9] = 0.31
$arr[0,10] = 0.31
$arr[0,11] = 0.31
$arr[0,12] = 0.31
$arr[0,13] = 0.31
$arr[0,14] = 0.31
$arr[0,15] = 0.31
$arr[0,16] = 0.31
$arr[0,17] = 0.31
$arr[0,18] = 0.31
$arr[0,19] = 0.31
$ws2.Range("AT2:BM2").Value = $arr

# ==== Sheet 3: strategy_id-5007 ====
$ws3 = $wb.Worksheets.Item(3)
$arr = New-Object 'object[,]' 1,4
$arr[0,0] = "variable_trajectory_group"
$arr[0,1] = "normalize_group"
$arr[0,2] = "trajgroup_no_vary_q"
$arr[0,3] = "uniform_scaling_q"
$ws3.Range("C1:F1").Value = $arr
$arr = New-Object 'object[,]' 1,2
$arr[0,0] = "max_55"
$arr[0,1] = "min_55"
$ws3.Range("H1:I1").Value = $arr
$ws3.Range("J1").Copy($ws3.Range("AT1:BM1"))
$arr = New-Object 'object[,]' 1,20
$arr[0,0] = 36
$arr[0,1] = 37
$arr[0,2] = 38
$arr[0,3] = 39
$arr[0,4] = 40
$arr[0,5] = 41
$arr[0,6] = 42
$arr[0,7] = 43
$arr[0,8] = 44
$arr[0,9] = 45
$arr[0,10] = 46
$arr[0,11] = 47
$arr[0,12] = 48
$arr[0,13] = 49
$arr[0,14] = 50
$arr[0,15] = 51
$arr[0,16] = 52
$arr[0,17] = 53
$arr[0,18] = 54
$arr[0,19] = 55
$ws3.Range("AT1:BM1").Value = $arr
$ws3.Range("C2").Value = 13
$arr = New-Object 'object[,]' 1,36
$arr[0,0] = 0.31
$arr[0,1] = 0.31
$arr[0,2] = 0.31
$arr[0,3] = 0.31
$arr[0,4] = 0.31
$arr[0,5] = 0.31
$arr[0,6] = 0.31
$arr[0,7] = 0.31
$arr[0,8] = 0.31
$arr[0,9] = 0.31
$arr[0,10] = 0.31
$arr[0,11] = 0.31
$arr[0,12] = 0.31
$arr[0,13] = 0.31
$arr[0,14] = 0.31
$arr[0,15] = 0.31
$arr[0,16] = 0.31
$arr[0,17] = 0.31
$arr[0,18] = 0.31
$arr[0,19] = 0.31
$arr[0,20] = 0.31
$arr[0,21] = 0.31
$arr[0,22] = 0.3100000000000001
$arr[0,23] = 0.31
$arr[0,24] = 0.31
$arr[0,25] = 0.3100000000000001
$arr[0,26] = 0.3099999999999999
$arr[0,27] = 0.31
$arr[0,28] = 0.31
$arr[0,29] = 0.3099999999999999
$arr[0,30] = 0.31
$arr[0,31] = 0.31
$arr[0,32] = 0.31
$arr[0,33] = 0.31
$arr[0,34] = 0.31
$arr[0,35] = 0.31
$ws3.Range("J2:AS2").Value = $arr
$arr = New-Object 'object[,]' 1,20
$arr[0,0] = 0.31
$arr[0,1] = 0.31
$arr[0,2] = 0.31
$arr[0,3] = 0.31
$arr[0,4] = 0.31
$arr[0,5] = 0.31
$arr[0,6] = 0.31
$arr[0,7] = 0.31
$arr[0,8] = 0.31
$arr[0,9] = 0.31
$arr[0,10] = 0.31
$arr[0,11] = 0.31
$arr[0,12] = 0.31
$arr[0,13] = 0.31
$arr[0,14] = 0.31
$arr[0,15] = 0.31
$arr[0,16] = 0.31
$arr[0,17] = 0.31
$arr[0,18] = 0.31
$arr[0,19] = 0.31
$ws3.Range("AT2:BM2").Value = $arr

# ==== Sheet 4: strategy_id-5009 ====
$ws4 = $wb.Worksheets.Item(4)
$arr = New-Object 'object[,]' 1,4
$arr[0,0] = "variable_trajectory_group"
$arr[0,1] = "normalize_group"
$arr[0,2] = "trajgroup_no_vary_q"
$arr[0,3] = "uniform_scaling_q"
$ws4.Range("C1:F1").Value = $arr
$arr = New-Object 'object[,]' 1,2
$arr[0,0] = "max_55"
$arr[0,1] = "min_55"
$ws4.Range("H1:I1").Value = $arr
$ws4.Range("J1").Copy($ws4.Range("AT1:BM1"))
$arr = New-Object 'object[,]' 1,20
$arr[0,0] = 36
$arr[0,1] = 37
$arr[0,2] = 38
$arr[0,3] = 39
$arr[0,4] = 40
$arr[0,5] = 41
$arr[0,6] = 42
$arr[0,7] = 43
$arr[0,8] = 44
$arr[0,9] = 45
$arr[0,10] = 46
$arr[0,11] = 47
$arr[0,12] = 48
$arr[0,13] = 49
$arr[0,14] = 50
$arr[0,15] = 51
$arr[0,16] = 52
$arr[0,17] = 53
$arr[0,18] = 54
$arr[0,19] = 55
$ws4.Range("AT1:BM1").Value = $arr
$ws4.Range("C2").Value = 13
$arr = New-Object 'object[,]' 1,36
$arr[0,0] = 0.31
$arr[0,1] = 0.31
$arr[0,2] = 0.31
$arr[0,3] = 0.31
$arr[0,4] = 0.31
$arr[0,5] = 0.31
$arr[0,6] = 0.31
$arr[0,7] = 0.31
$arr[0,8] = 0.31
$arr[0,9] = 0.31
$arr[0,10] = 0.31
$arr[0,11] = 0.31
$arr[0,12] = 0.31
$arr[0,13] = 0.31
$arr[0,14] = 0.31
$arr[0,15] = 0.31
$arr[0,16] = 0.31
$arr[0,17] = 0.31
$arr[0,18] = 0.31
$arr[0,19] = 0.31
$arr[0,20] = 0.31
$arr[0,21] = 0.31
$arr[0,22] = 0.3100000000000001
$arr[0,23] = 0.31
$arr[0,24] = 0.31
$arr[0,25] = 0.3100000000000001
$arr[0,26] = 0.3099999999999999
$arr[0,27] = 0.31
$arr[0,28] = 0.31
$arr[0,29] = 0.3099999999999999
$arr[0,30] = 0.31
$arr[0,31] = 0.31
$arr[0,32] = 0.31
$arr[0,33] = 0.31
$arr[0,34] = 0.31
$arr[0,35] = 0.31
$ws4.Range("J2:AS2").Value = $arr
$arr = New-Object 'object[,]' 1,20
$arr[0,0] = 0.31
$arr[0,1] = 0.31
$arr[0,2] = 0.31
$arr[0,3] = 0.31
$arr[0,4] = 0.31
$arr[0,5] = 0.31
$arr[0,6] = 0.31
$arr[0,7] = 0.31
$arr[0,8] = 0.31
$arr[0,9] = 0.31
$arr[0,10] = 0.31
$arr[0,11] = 0.31
$arr[0,12] = 0.31
$arr[0,13] = 0.31
$arr[0,14] = 0.31
$arr[0,15] = 0.31
$arr[0,16] = 0.31
$arr[0,17] = 0.31
$arr[0,18] = 0.31
$arr[0,19] = 0.31
$ws4.Range("AT2:BM2").Value = $arr
